# Update poisson.xlsx (Sheet1) for "semana 30 de 2024"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C3" = 3;   "D3" = 3;   "E3" = 0.22;
    "D4" = 0;   "E4" = 1;
    "C5" = 3;   "D5" = 8;   "E5" = 0.01;
    "C6" = 0;   "D6" = 108;
    "D7" = 1;   "E7" = 0.27;
    "C11" = 48; "D11" = 33;
    "C13" = 1;  "D13" = 3;  "E13" = 0.06;
    "C14" = 3;  "D14" = 7;  "E14" = 0.02;
    "C15" = 29;
    "C16" = 4;  "E16" = 0.02;
    "C17" = 1;  "E17" = 0.37;
    "D18" = 1;  "E18" = 0.27;
    "C19" = 10; "D19" = 4;  "E19" = 0.02;
    "C20" = 2;  "E20" = 0.14;
    "C21" = 5;  "D21" = 1;  "E21" = 0.03;
    "D25" = 1;
    "C27" = 4;  "D27" = 6;  "E27" = 0.1;
    "D28" = 1;
    "C30" = 1;  "E30" = 0.37;
    "D31" = 1;  "E31" = 0;
    "C32" = 1;  "E32" = 0.18;
    "D33" = 5;  "E33" = 0.13;
    "C34" = 7;  "D34" = 4;  "E34" = 0.09;
    "C35" = 10; "D35" = 7;  "E35" = 0.09;
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
